# Update countries & provincias Spain
# Applies updated COVID case numbers for several countries and reflects the
# resulting re-sort (by total cases, column B) for the affected blocks of
# rows, plus refreshes the "last updated" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh "last updated" timestamp (A1) ------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 12:52"

# --- Straightforward numeric refreshes (row/country unchanged) ---------

# Row 8: Alemania
$ws.Range("B8").Value = 148587
$ws.Range("C8").Value = 134
$ws.Range("D8").Value = 99400
$ws.Range("E8").Value = 44097
$ws.Range("F8").Value = 2908
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 5090

# Row 11: Iran
$ws.Range("B11").Value = 85996
$ws.Range("C11").Value = 1194
$ws.Range("D11").Value = 63113
$ws.Range("E11").Value = 17492
$ws.Range("F11").Value = 3311
$ws.Range("G11").Value = 94
$ws.Range("H11").Value = 5391

# Row 68: Uzbekistan
$ws.Range("D68").Value = 407
$ws.Range("E68").Value = 1278
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 7

# --- Re-sorted blocks: values updated and countries shift rows ---------

# Rows 81-83 (Ghana moves above Afganistan & Cuba)
$ws.Range("A81").Value = "Ghana"
$ws.Range("B81").Value = 1154
$ws.Range("C81").Value = 112
$ws.Range("D81").Value = 99
$ws.Range("E81").Value = 1046
$ws.Range("F81").Value = 4
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 9

$ws.Range("A82").Value = "Afganistan"
$ws.Range("B82").Value = 1143
$ws.Range("C82").Value = 51
$ws.Range("D82").Value = 166
$ws.Range("E82").Value = 937
$ws.Range("F82").Value = 7
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 40

$ws.Range("A83").Value = "Cuba"
$ws.Range("B83").Value = 1137
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 309
$ws.Range("E83").Value = 790
$ws.Range("F83").Value = 18
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 38

# Rows 107-109 (Senegal moves above Jordania & Taiwan)
$ws.Range("A107").Value = "Senegal"
$ws.Range("B107").Value = 442
$ws.Range("C107").Value = 30
$ws.Range("D107").Value = 253
$ws.Range("E107").Value = 183
$ws.Range("F107").Value = 1
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 6

$ws.Range("A108").Value = "Jordania"
$ws.Range("B108").Value = 428
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 297
$ws.Range("E108").Value = 124
$ws.Range("F108").Value = 5
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 7

$ws.Range("A109").Value = "Taiwan"
$ws.Range("B109").Value = 426
$ws.Range("C109").Value = 1
$ws.Range("D109").Value = 236
$ws.Range("E109").Value = 184
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 6

# Rows 165-166 (Nepal moves above Macao)
$ws.Range("A165").Value = "Nepal"
$ws.Range("B165").Value = 45
$ws.Range("C165").Value = 3
$ws.Range("D165").Value = 7
$ws.Range("E165").Value = 38
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 0

$ws.Range("A166").Value = "Macao"
$ws.Range("B166").Value = 45
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 24
$ws.Range("E166").Value = 21
$ws.Range("F166").Value = 1
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0
